# Update the NPV figures in the economic parameters table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$1812 Billion'
$ws.Range("B3").Value = '$1552 Billion'
$ws.Range("B4").Value = '$1421 Billion'

# Match the author's new active-cell selection on save.
$ws.Range("B5").Select()
